$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells holding numeric-looking text that must stay text (not auto-converted to numbers)
$numericTextCells = @("B4","D4","F4","H4","B5","D5","F5","H5","B6","D6")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 4: was BSKgezien=Yes -> now CorrCa24u (with new coef/CI/Chi values)
$ws.Range("A4").Value = "CorrCa24u"
$ws.Range("B4").Value = "1.458"
$ws.Range("C4").Value = "[0.112; 17.271]"
$ws.Range("D4").Value = "8.6"
$ws.Range("F4").Value = "1.426"
$ws.Range("G4").Value = "[0.121; 15.428]"
$ws.Range("H4").Value = "8.2"

# Row 5: was CorrCa24u -> now BSKgezien - No:Yes (with new coef/CI/Chi values)
$ws.Range("A5").Value = "BSKgezien - No:Yes"
$ws.Range("B5").Value = "3.479"
$ws.Range("C5").Value = "[2.906; 3.042]"
$ws.Range("D5").Value = "7.2"
$ws.Range("F5").Value = "3.778"
$ws.Range("G5").Value = "[1.335; 7.648]"
$ws.Range("H5").Value = "8.9"

# Row 6: Age_Years keeps its name, but coef/CI values change
$ws.Range("A6").Value = "Age_Years"
$ws.Range("B6").Value = "1.157"
$ws.Range("C6").Value = "[0.456; 2.827]"
$ws.Range("D6").Value = "0.2"

# Row 7: Sex=Male -> Sex - Male:Female
$ws.Range("A7").Value = "Sex - Male:Female"

# Row 8: surgery_type=total -> surgery_type - completion:total
$ws.Range("A8").Value = "surgery_type - completion:total"

# Row 9: CHKD=Yes -> CHKD - Yes:No
$ws.Range("A9").Value = "CHKD - Yes:No"

# Restore the original (default) cell formatting now that text values are safely stored
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).ClearFormats()
}
